# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" text with the new rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.57 = 5713.84 pesos`n✅ 5713.84 pesos = 1.56 = 936.69 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- tasas: update the N10/O10 and N12/O12 rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 637.4
$wsTasas.Range("O10").Value = 3642
$wsTasas.Range("N12").Value = 3660
$wsTasas.Range("O12").Value = 600
